$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FERNANDEZ VALDERAS ERNESTO ALI"
$ws.Range("B2").Value = 157

$ws.Range("A3").Value = "GUTIERREZ CARLOS TERESA DE JESUS"
$ws.Range("B3").Value = 156

$ws.Range("A4").Value = "VALLE MAGALLAN EDUAR"
$ws.Range("B4").Value = 138

$ws.Range("A5").Value = "CONTRERAS VALDERRAMA JULIA ALEJANDRA"
$ws.Range("B5").Value = 136

$ws.Range("B6").Value = 109

$ws.Range("A7").Value = "CAMACHO LINARES JUDITH ARLETT"
$ws.Range("B7").Value = 107

$ws.Range("A8").Value = "SEVERINO AVALOS MARJORIE ISABEL"
$ws.Range("B8").Value = 106

$ws.Range("B9").Value = 101

$ws.Range("A10").Value = "ROMERO CHANAME YOSSELY TRINIDAD"
$ws.Range("B10").Value = 100

$ws.Range("A12").Value = "ZEVALLOS PACHECO ZOILA XIMENA"
$ws.Range("B12").Value = 95

$ws.Range("A13").Value = "BALLENA ESQUÉN ASTRID CAROLINA"
$ws.Range("B13").Value = 93

$ws.Range("A14").Value = "SENADOR ARBOLEDA GIANCARLOS EXEBIO"
$ws.Range("B14").Value = 91
